$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new weekly rows before the current row 168 (shifts old 168-240
# down to 170-242, growing the used range from A1:R240 to A1:R242).
$ws.Rows("168:169").Insert()

# New "Primera" quality row for Vega Monumental Concepción - Betarraga,
# week of 2021-12-23 (serial 44553).
$ws.Cells.Item(168, 1).Value = 11
$ws.Cells.Item(168, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(168, 3).Value = "Bíobío"
$ws.Cells.Item(168, 4).Value = 44553
$ws.Cells.Item(168, 5).Value = 8
$ws.Cells.Item(168, 6).Value = 100114014
$ws.Cells.Item(168, 7).Value = "Betarraga"
$ws.Cells.Item(168, 8).Value = "Sin especificar"
$ws.Cells.Item(168, 9).Value = "Primera"
$ws.Cells.Item(168, 10).Value = 600
$ws.Cells.Item(168, 11).Value = 600
$ws.Cells.Item(168, 12).Value = 700
$ws.Cells.Item(168, 13).Value = 650
$ws.Cells.Item(168, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(168, 15).Value = "Región Metropolitana"
$ws.Cells.Item(168, 16).Value = 130
$ws.Cells.Item(168, 17).Value = 5
$ws.Cells.Item(168, 18).Value = "Hortaliza"

# Matching "Segunda" quality row for the same market/date.
$ws.Cells.Item(169, 1).Value = 11
$ws.Cells.Item(169, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(169, 3).Value = "Bíobío"
$ws.Cells.Item(169, 4).Value = 44553
$ws.Cells.Item(169, 5).Value = 8
$ws.Cells.Item(169, 6).Value = 100114014
$ws.Cells.Item(169, 7).Value = "Betarraga"
$ws.Cells.Item(169, 8).Value = "Sin especificar"
$ws.Cells.Item(169, 9).Value = "Segunda"
$ws.Cells.Item(169, 10).Value = 300
$ws.Cells.Item(169, 11).Value = 500
$ws.Cells.Item(169, 12).Value = 500
$ws.Cells.Item(169, 13).Value = 500
$ws.Cells.Item(169, 14).Value = '$/paquete 5 unidades'
$ws.Cells.Item(169, 15).Value = "Región Metropolitana"
$ws.Cells.Item(169, 16).Value = 100
$ws.Cells.Item(169, 17).Value = 5
$ws.Cells.Item(169, 18).Value = "Hortaliza"
